$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (5th column), shifting Weapon.. etc right by one
$ws.Columns("E").Insert()

# Update C2 (character level) from numeric 90 to the text "90+"
$ws.Range("C2").Value = "90+"

# Set header for new column E1
$ws.Range("E1").Value = "技能等级"

# Set data for new column E2
$ws.Range("E2").Value = "9,9,10"

# Fill in the new data columns J2 through AP2 as per the skill reader data
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 311
$ws.Range("L2").Value = 46.6
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 47800
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 65
$ws.Range("U2").Value = 90
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 51.8
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0
$ws.Range("AO2").Value = 0
$ws.Range("AP2").Value = 0

# New last column (AP) gets its own width, distinct from AO's bestFit width
$ws.Range("AP1").ColumnWidth = 14

# Update selection to match the target file (E2 selected)
$ws.Range("E2").Select()
